# Fix typo and naming description in columns of template sheets (#541)
#
# - "QBIC sample ids*" -> "QBIC sample id*"   (Metadata + Allowed-Values header row)
# - "QBIC sample ids"  -> "QBIC sample id"    (Property information description table)
# - Long winded description of the QBIC sample id column shortened
# - The "single-end" / "paired-end" example values in the Allowed-Values sheet
#   are replaced by a generic "Free text" entry (single-end) and removed (paired-end)

$wb = $excel.ActiveWorkbook

$metadata = $wb.Sheets.Item("Metadata")
$propertyInfo = $wb.Sheets.Item("Property information")
$allowedValues = $wb.Sheets.Item("Allowed-Values")

# Header typo: "QBIC sample ids*" -> "QBIC sample id*"
$metadata.Range("A1").Value = "QBIC sample id*"
$allowedValues.Range("A1").Value = "QBIC sample id*"

# Property information description table
$propertyInfo.Range("A2").Value = "QBIC sample id"
$propertyInfo.Range("C2").Value = "Each measurement need to be linked to at least on analyte sample."

# Allowed-Values sheet: drop the single-end / paired-end example values
$allowedValues.Range("F2").Value = "Free text"
$allowedValues.Range("F3").ClearContents()
